# Applies the crypto price/volume refresh described in the commit diff.
# Column D price strings that would otherwise be auto-parsed as numbers are
# written with a leading text quote-prefix (') so Excel keeps them stored as
# text, matching the source workbook where every Price/Volume cell is a string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.264.94'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '3.626.85'
$ws.Range('E3').Value = '  +3.61%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''604.32'
$ws.Range('E5').Value = '  +0.97%  '
$ws.Range('D6').Value = '''196.30'
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('E7').Value = '  +0.95%  '
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('E10').Value = '  +0.16%  '
$ws.Range('D11').Value = '''54.08'
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('D12').Value = '''0.0000307'
$ws.Range('E12').Value = '  +2.20%  '
$ws.Range('E13').Value = '  +0.62%  '
$ws.Range('D14').Value = '4.201.97'
$ws.Range('E14').Value = '  +3.61%  '
$ws.Range('D15').Value = '''13.27'
$ws.Range('E15').Value = '  +5.70%  '
$ws.Range('D16').Value = '''595.04'
$ws.Range('E16').Value = '  -1.85%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.641.48'
$ws.Range('E17').Value = '  +4.04%  '
$ws.Range('D18').Value = '70.452.21'
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '''19.24'
$ws.Range('E19').Value = '  +1.38%  '
$ws.Range('E20').Value = '  +1.57%  '
$ws.Range('D21').Value = '''0.999'
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('E22').Value = '  -2.10%  '
$ws.Range('E23').Value = '  +1.33%  '
$ws.Range('D24').Value = '''102.39'
$ws.Range('E24').Value = '  -1.71%  '
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('D26').Value = '''3.04'
$ws.Range('E26').Value = '  -0.55%  '
$ws.Range('D27').Value = '''10.83'
$ws.Range('E27').Value = '  -1.06%  '
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('E29').Value = '  +1.94%  '
$ws.Range('D30').Value = '''4.77'
$ws.Range('E30').Value = '  +6.72%  '
$ws.Range('D31').Value = '''7.21'
$ws.Range('E31').Value = '  +1.64%  '
$ws.Range('E32').Value = '  -2.13%  '
$ws.Range('E33').Value = '  +2.57%  '
$ws.Range('D34').Value = '0.0₃0911'
$ws.Range('E34').Value = '  +12.78%  '
$ws.Range('D35').Value = '''63.28'
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('D36').Value = '3.923.75'
$ws.Range('E36').Value = '  +5.34%  '
$ws.Range('D37').Value = '''3.15'
$ws.Range('E37').Value = '  +3.64%  '
$ws.Range('D38').Value = '''529.82'
$ws.Range('E38').Value = '  +5.78%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').Value = '''37.55'
$ws.Range('E40').Value = '  +2.45%  '
$ws.Range('D41').Value = '''0.393'
$ws.Range('E41').Value = '  +0.63%  '
$ws.Range('D42').Value = '''3.55'
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('E43').Value = '  -1.04%  '
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').Value = '''2.87'
$ws.Range('E45').Value = '  +2.24%  '
$ws.Range('D46').Value = '''3.36'
$ws.Range('E46').Value = '  +1.02%  '
$ws.Range('E47').Value = '  +0.81%  '
$ws.Range('E48').Value = '  -1.02%  '
$ws.Range('D50').Value = '''0.000255'
$ws.Range('E50').Value = '  +5.93%  '
$ws.Range('D51').Value = '''1.32'
$ws.Range('E51').Value = '  +4.02%  '
